$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Remove the duplicate "Contact" row (row 11) - rows 12-15 shift up to 11-14
$ws.Rows.Item(11).Delete()

# Version: 5.0.0 -> 6.0.0
$ws.Cells.Item(3, 2).Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value: (empty) -> Alvearie Team
$ws.Cells.Item(9, 2).Value = "Alvearie Team"

# Former "Contact" / "No display for ContactDetail" row becomes "Jurisdiction" / "United States of America"
$ws.Cells.Item(10, 1).Value = "Jurisdiction"
$ws.Cells.Item(10, 2).Value = "United States of America"
